$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("E13").Value = 'Dict mit type : {Votetypen}, dummy : {"True", "False"} (ob der Client mitvoten soll) und players : Liste der Spielernamen '

$ws.Range("E14").Select()
